# Apply "Updated codebook and data dictionary" changes to KIM_codebook.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the previously-empty frequency/percent/missing counts for the
#    "Obesity" (row 44/45) variable.
# ---------------------------------------------------------------------------
$ws.Range("J44").Value = 8657
$ws.Range("K44").Value = 63.56
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("M44").NumberFormat = "0.00"

$ws.Range("J45").Value = 4963
$ws.Range("K45").Value = 36.44

# ---------------------------------------------------------------------------
# 2) Renumber the "Number" column for the variables that now sit one row
#    further down the codebook (28->29, 29->30, ... 34->35) because a new
#    "Underweight" variable was inserted ahead of them logically.
# ---------------------------------------------------------------------------
$ws.Range("A46").Value = 29
$ws.Range("A48").Value = 30
$ws.Range("A50").Value = 31
$ws.Range("A51").Value = 32
$ws.Range("A52").Value = 33
$ws.Range("A53").Value = 34
$ws.Range("A54").Value = 35

# ---------------------------------------------------------------------------
# 3) Append the two new variables to the codebook: "Underweight" (#36) and
#    "Agegroup" (#37), each with their own frequency tables.
# ---------------------------------------------------------------------------

# --- Variable 36: Underweight -------------------------------------------------
$ws.Range("A58").Value = 36
$ws.Range("B58").Value = "Underweight"
$ws.Range("C58").Value = "Underweight (BMI <18.5 kg/m^2)"
$ws.Range("D58").Value = "Character"

$ws.Range("H58").Value = 0
$ws.Range("H58").HorizontalAlignment = -4152
$ws.Range("I58").Value = "No"
$ws.Range("J58").Value = 13067
$ws.Range("K58").Value = 95.94
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("M58").NumberFormat = "0.00"

$ws.Range("H59").Value = 1
$ws.Range("H59").HorizontalAlignment = -4152
$ws.Range("I59").Value = "Yes"
$ws.Range("J59").Value = 553
$ws.Range("K59").Value = 4.0599999999999996

# --- Variable 37: Agegroup ----------------------------------------------------
$ws.Range("A60").Value = 37
$ws.Range("B60").Value = "Agegroup"
$ws.Range("C60").Value = "Age group (years)"
$ws.Range("D60").Value = "Character"

$ws.Range("H60").Value = 1
$ws.Range("H60").HorizontalAlignment = -4152
$ws.Range("I60").Value = "20-29 years"
$ws.Range("J60").Value = 1054
$ws.Range("L60").Value = 33
$ws.Range("M60").Value = 0.24
$ws.Range("M60").NumberFormat = "0.00"

$ws.Range("H61").Value = 2
$ws.Range("H61").HorizontalAlignment = -4152
$ws.Range("I61").Value = "30-39 years"
$ws.Range("J61").Value = 3179

$ws.Range("H62").Value = 3
$ws.Range("H62").HorizontalAlignment = -4152
$ws.Range("I62").Value = "40-49 years"
$ws.Range("J62").Value = 3136

$ws.Range("H63").Value = 4
$ws.Range("H63").HorizontalAlignment = -4152
$ws.Range("I63").Value = "50-59 years"
$ws.Range("J63").Value = 3574

$ws.Range("H64").Value = 5
$ws.Range("H64").HorizontalAlignment = -4152
$ws.Range("I64").Value = "60-69 years"
$ws.Range("J64").Value = 1936

$ws.Range("H65").Value = 6
$ws.Range("H65").HorizontalAlignment = -4152
$ws.Range("I65").Value = "70 years and older"
$ws.Range("J65").Value = 708

# K60 gets its own formula, K61:K65 share one formula (mirrors the shared
# formula Excel created for this column in the authored workbook).
$ws.Range("K60").Formula = "=J60/13620*100"
$ws.Range("K60").NumberFormat = "0.00"

$ws.Range("K61:K65").Formula = "=J61/13620*100"
$ws.Range("K61:K65").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 4) Update the view state to match where the author left the cursor.
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("L61").Select()
